$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append newly collected tweets (sharonjcheboi) to column A, one per
# odd-numbered row gap, matching the pasted layout from the source sheet.
$ws.Cells.Item(278, 1).Value2 = 'While working to minimise impact of #Covid19 on food security, we need to document what worked, what didn''t, what could have been done different or better or what should have been proporitised more; globally to individual level. #CovidFoodFuture  @sharonjcheboi'
$ws.Cells.Item(280, 1).Value2 = '#Kenya Justifiably, the government''s focus has shifted to managing #Covid19. Focus should also be drawn towards food security, which is just as important as health security. Poor nutrition promotes susceptibility of an individual to infection. #CovidFoodFuture  @sharonjcheboi'
$ws.Cells.Item(282, 1).Value2 = '@FAO warned of a massive increase in locusts across East Africa, threatening the already food insecure nations. Control has been slow as borders close due to #Covid19, hence hindering pesticide import.#Covid19FoodFuture @sharonjcheboi'
$ws.Cells.Item(284, 1).Value2 = '#Kenya #Covid19 expresses the importance of self reliance on food security of a country or household. Encouraging self reliance for those who can access urban agricultural spaces is important. Eg kale, spinach, managu and herbs, potatoes, carrots #Covid19FoodFuture @sharonjcheboi'
$ws.Cells.Item(286, 1).Value2 = 'To minimise the impact of #Covid19 on food security in #Africa, What actions should be taken to minimise the impact on food distribution? #Covid19FoodFuture @sharonjcheboi'
$ws.Cells.Item(288, 1).Value2 = 'A crisis such as #Covid19 might be a health crisis. But has the possibility of causing a food crisis. Food systems have been affected and a food crisis within a health crisis should be prevented. Here''s how, by @WorldBank https://bit.ly/2K2DPGV  @sharonjcheboi #Covid19FoodFuture'
$ws.Cells.Item(290, 1).Value2 = 'This article illustrates the impact of #Covid19 on food joints: Due to stock purchase of products such as tea and coffee, consumers have resorted to home consumption of beverages they would usually consume in hotels. https://bit.ly/2V6c8TP @sharonjcheboi #Covid19FoodFuture'
$ws.Cells.Item(292, 1).Value2 = 'Farmers globally are facing poor market access as per the poll. Higher value crops which require higher labour are lost before and after harvest, reducing profits. Due to income losses, fewer consumers maybe able to affords these expensive crops. #Covid19FoodFuture  @sharonjcheboi'
$ws.Cells.Item(294, 1).Value2 = 'Kenya, through @FredMatiangi recently banned uncoordinated food distribution within communities. Was it a good move? Has it provided enough measures to cushion its most vulnerable? Does it disregard food as a basic need? What are your thoughts? #CovidFoodFuture @sharonjcheboi'
$ws.Cells.Item(296, 1).Value2 = 'In Kenya and all over Africa, many farmers have lost their bulk buyers as schools and markets shut down. Aside from platforms such as @mkulimayoung. How else can farmers cushion themselves and market their produce to minimise complete loss? #CovidFoodFuture @sharonjcheboi'
$ws.Cells.Item(298, 1).Value2 = 'Yesterday in Kibera, a failed food distribution strategy left people injured and 2 women dead during the chaos. During #covid19, cash transfers to caregivers would be more dignified and adheres to social distancing rules by @WHO #Covid19FoodFuture @sharonjcheboi'
$ws.Cells.Item(300, 1).Value2 = 'This article illustrates the impact of #Covid19 on food joints. Due to stock purchase of products such as tea and coffee, consumers have resorted to come consumption of beverages they would usually consume in hotels. #Covid19FoodFuture @sharonjcheboi bit.ly/2V6c8TP'
$ws.Cells.Item(302, 1).Value2 = 'Farmers globally are facing poor market access as per the poll. Higher value crops which require high labour are lost before and after harvest, reducing profits. Due to income losses, fewer consumers may be able to afford these crops. #Covid19FoodFuture @sharonjcheboi'
$ws.Cells.Item(306, 1).Value2 = 'In Kenya, @wawiranjiru of @Food4Education has set towards providing 6000 households with food staples during #covid19 pandemic. Read more below: https://africasustainabilitymatters.com/social-entrepreneur-28-to-feed-poor-kenyan-homes-amid-coronavirus/amp/?__twitter_impression=true @sharonjcheboi'
$ws.Cells.Item(308, 1).Value2 = 'In Kenya the @MOH_Kenya announced the first confirmed #covid19 case in Uasin Gishu County. A local farmer has expressed their challenge of accessing labour and inputs due to peoples'' fear of infection and movement restrictions. #covid19foodfuture @sharonjcheboi'
$ws.Cells.Item(310, 1).Value2 = 'Farm inputs distribution has been restricted and farm input costs have gone higher as a result of higher transportation costs. #Covid19food19future @sharonjcheboi'
$ws.Cells.Item(312, 1).Value2 = 'Kenyan farmers lost 30 - 35% of 2019 harvest as a result of heavy off season rains. The shortage will very soon be felt in the market. Higher prices for a population that has lost its livelihoods will increase preexistent food insecurity. @sharonjcheboi  #Covid19FoodFuture'
$ws.Cells.Item(314, 1).Value2 = 'Most food security measures in Kenya focus on production while accessibility is equally important. Due to #COVID19, markets have declined causing massive post-harvest losses. @sharonjcheboi'
$ws.Cells.Item(316, 1).Value2 = 'Following the heavy rains in January, the recurring locust invasion in February and the #COVID19 pandemic in March, the already vulnerable Kenyan #foodsystem will suffer increasing food insecurity There''s therefore a need to developed resilient systems. @sharonjcheboi'

# Restore the selection/scroll position left by the editor.
$ws.Range("A278:A318").Select()
